# Insert three new paragraphs at the very top of the document body,
# before the existing first paragraph:
#   1) "I am bob, every day I need to write a report, "
#   2) "Feel so boring1"
#   3) an empty paragraph (no run at all)
#
# The rest of the document is left untouched.

$d = $word.ActiveDocument

$first = $d.Paragraphs(1)

# Create three blank paragraphs ahead of the current first paragraph.
# Each call inserts a new empty paragraph immediately before $first,
# so after three calls, paragraphs 1, 2, 3 are the new (still empty)
# paragraphs and paragraph 4 is the original first paragraph.
$first.Range.InsertParagraphBefore()
$first.Range.InsertParagraphBefore()
$first.Range.InsertParagraphBefore()

# Fill in the text for the first two new paragraphs.
$d.Paragraphs(1).Range.Text = "I am bob, every day I need to write a report, "
$d.Paragraphs(2).Range.Text = "Feel so boring1"

# The third new paragraph must stay completely empty - i.e. contain no
# run/text element whatsoever, matching a paragraph that was created
# but never typed into. Assigning an empty string to Range.Text leaves
# behind a run with an empty <w:t>, so instead give it one placeholder
# character and then delete that single-character range: deleting an
# actual (non-zero-length) character range removes the run entirely,
# leaving just the bare paragraph.
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "X"
$p3 = $d.Paragraphs(3)
$charRange = $d.Range($p3.Range.Start, $p3.Range.Start + 1)
$charRange.Delete()
